$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("D2").Value = '62.634.16'
    $ws.Range("E2").Value = '  +1.14%  '
    # Row 3
    $ws.Range("D3").Value = '3.460.57'
    $ws.Range("E3").Value = '  +1.04%  '
    # Row 4
    $ws.Range("E4").Value = '  +0.01%  '
    # Row 5
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = '413.43'
    $ws.Range("E5").Value = '  +1.12%  '
    # Row 6
    $ws.Range("D6").NumberFormat = "@"
    $ws.Range("D6").Value = '130.10'
    $ws.Range("E6").Value = '  +1.27%  '
    # Row 7
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = '0.621'
    $ws.Range("E7").Value = '  -1.73%  '
    # Row 8
    $ws.Range("E8").Value = '  +0.04%  '
    # Row 9
    $ws.Range("D9").NumberFormat = "@"
    $ws.Range("D9").Value = '0.723'
    $ws.Range("E9").Value = '  -2.04%  '
    # Row 10
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = '0.143'
    $ws.Range("E10").Value = '  +1.09%  '
    # Row 11
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = '42.48'
    $ws.Range("E11").Value = '  -0.50%  '
    # Row 12
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = '9.64'
    $ws.Range("E12").Value = '  +5.49%  '
    # Row 13
    $ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
    $ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    $ws.Range("D13").Value = '4.016.97'
    $ws.Range("E13").Value = '  +1.36%  '
    # Row 14
    $ws.Range("B14").Value = 'ShibaInu'
    $ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = '0.0000216'
    $ws.Range("E14").Value = '  -1.86%  '
    # Row 15
    $ws.Range("E15").Value = '  -0.26%  '
    # Row 16
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = '20.45'
    $ws.Range("E16").Value = '  -4.50%  '
    # Row 17
    $ws.Range("D17").Value = '3.482.85'
    $ws.Range("E17").Value = '  +3.42%  '
    # Row 18
    $ws.Range("E18").Value = '  +1.58%  '
    # Row 19
    $ws.Range("E19").Value = '  -1.10%  '
    # Row 20
    $ws.Range("D20").Value = '62.629.64'
    $ws.Range("E20").Value = '  +1.13%  '
    # Row 21
    $ws.Range("D21").NumberFormat = "@"
    $ws.Range("D21").Value = '466.29'
    $ws.Range("E21").Value = '  +3.03%  '
    # Row 22
    $ws.Range("D22").NumberFormat = "@"
    $ws.Range("D22").Value = '90.59'
    $ws.Range("E22").Value = '  -1.64%  '
    # Row 23
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = '3.26'
    $ws.Range("E23").Value = '  +1.68%  '
    # Row 24
    $ws.Range("D24").NumberFormat = "@"
    $ws.Range("D24").Value = '13.32'
    $ws.Range("E24").Value = '  +2.56%  '
    # Row 25
    $ws.Range("D25").NumberFormat = "@"
    $ws.Range("D25").Value = '10.60'
    $ws.Range("E25").Value = '  +20.64%  '
    # Row 26
    $ws.Range("D26").NumberFormat = "@"
    $ws.Range("D26").Value = '3.29'
    $ws.Range("E26").Value = '  +1.77%  '
    # Row 27
    $ws.Range("D27").NumberFormat = "@"
    $ws.Range("D27").Value = '33.24'
    $ws.Range("E27").Value = '  +0.98%  '
    # Row 28
    $ws.Range("E28").Value = '  +0.42%  '
    # Row 29
    $ws.Range("E29").Value = '  -1.09%  '
    # Row 30
    $ws.Range("D30").NumberFormat = "@"
    $ws.Range("D30").Value = '11.96'
    $ws.Range("E30").Value = '  -0.25%  '
    # Row 31
    $ws.Range("D31").NumberFormat = "@"
    $ws.Range("D31").Value = '2.64'
    $ws.Range("E31").Value = '  -2.72%  '
    # Row 32
    $ws.Range("E32").Value = '  -2.42%  '
    # Row 33
    $ws.Range("E33").Value = '  -1.99%  '
    # Row 34
    $ws.Range("D34").NumberFormat = "@"
    $ws.Range("D34").Value = '40.70'
    $ws.Range("E34").Value = '  -5.30%  '
    # Row 35
    $ws.Range("E35").Value = '  +0.15%  '
    # Row 36
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = '58.46'
    $ws.Range("E36").Value = '  +7.71%  '
    # Row 37
    $ws.Range("D37").NumberFormat = "@"
    $ws.Range("D37").Value = '0.0488'
    $ws.Range("E37").Value = '  -2.72%  '
    # Row 38
    $ws.Range("E38").Value = '  +0.13%  '
    # Row 39
    $ws.Range("D39").NumberFormat = "@"
    $ws.Range("D39").Value = '3.06'
    $ws.Range("E39").Value = '  +4.76%  '
    # Row 40
    $ws.Range("D40").NumberFormat = "@"
    $ws.Range("D40").Value = '3.34'
    $ws.Range("E40").Value = '  -1.18%  '
    # Row 41
    $ws.Range("B41").Value = 'WEMIXToken'
    $ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    $ws.Range("D41").NumberFormat = "@"
    $ws.Range("D41").Value = '2.70'
    $ws.Range("E41").Value = '  +7.19%  '
    # Row 42
    $ws.Range("B42").Value = 'TheGraph'
    $ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    $ws.Range("D42").NumberFormat = "@"
    $ws.Range("D42").Value = '0.320'
    # Row 43
    $ws.Range("B43").Value = 'Stellar'
    $ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = '0.133'
    $ws.Range("E43").Value = '  -1.08%  '
    # Row 44
    $ws.Range("D44").NumberFormat = "@"
    $ws.Range("D44").Value = '146.00'
    $ws.Range("E44").Value = '  +2.82%  '
    # Row 45
    $ws.Range("B45").Value = 'ARBITRUM'
    $ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    $ws.Range("D45").NumberFormat = "@"
    $ws.Range("D45").Value = '2.07'
    $ws.Range("E45").Value = '  +4.01%  '
    # Row 46
    $ws.Range("B46").Value = 'NEARProtocol'
    $ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    $ws.Range("D46").NumberFormat = "@"
    $ws.Range("D46").Value = '4.34'
    $ws.Range("E46").Value = '  +1.97%  '
    # Row 47
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = '2.41'
    $ws.Range("E47").Value = '  +11.77%  '
    # Row 48
    $ws.Range("D48").Value = '0.0₃0554'
    $ws.Range("E48").Value = '  +31.15%  '
    # Row 49
    $ws.Range("D49").NumberFormat = "@"
    $ws.Range("D49").Value = '16.33'
    $ws.Range("E49").Value = '  -2.12%  '
    # Row 50
    $ws.Range("D50").NumberFormat = "@"
    $ws.Range("D50").Value = '22.17'
    $ws.Range("E50").Value = '  -0.89%  '
    # Row 51
    $ws.Range("E51").Value = '  +0.89%  '